# Data Source.xlsx update
# - Sheet1 ("Stage I growth rate source"): add a Source note for the CPI Medicare row
# - Sheet2 ("Stage II base year source"): insert a new "Column #" reference column
#   (between the existing "Source table" and "Source description" columns) and
#   populate it, then update which sheet/cell is active in the saved view.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Stage I growth rate source")
$ws2 = $wb.Worksheets.Item("Stage II base year source")

# --- Sheet1: new "Source" note for CPI Medicare (row 17) ---
$ws1.Range("E17").Value = "CPI Medicare"

# --- Sheet2: insert a new column at E, shifting old E:F to F:G ---
$ws2.Columns.Item(5).Insert()

# Fill in the new column E with the underlying source-table column references
$ws2.Range("D8").Value  = "Table 2.3"
$ws2.Range("D15").Value = "Table 1.4"
$ws2.Range("D16").Value = "Table 1.4"

$ws2.Range("E4").Value  = "Column 66"
$ws2.Range("E5").Value  = "Columns 14 + 27 + 53"
$ws2.Range("E6").Value  = "Column 40"
$ws2.Range("E7").Value  = "Column 69"
$ws2.Range("E8").Value  = "Column 6"
$ws2.Range("E9").Value  = "Column 8"
$ws2.Range("E10").Value = "Column 12"
$ws2.Range("E11").Value = "Column 20"
$ws2.Range("E12").Value = "Column 22"
$ws2.Range("E13").Value = "Column 26"
$ws2.Range("E14").Value = "Column 38"
$ws2.Range("E15").Value = "Columns 52 + 56 + 60"
$ws2.Range("E16").Value = "Columns 54 + 58 + 62"
$ws2.Range("E17").Value = "Column 70"
$ws2.Range("E18").Value = "Column 68"
$ws2.Range("E19").Value = "Column 6"
$ws2.Range("E20").Value = "Column 6"
$ws2.Range("E21").Value = "Column 6"
$ws2.Range("E22").Value = "Column 6"
$ws2.Range("E23").Value = "Column 6"
$ws2.Range("E24").Value = "Column 6"
$ws2.Range("E25").Value = "Column 6"
$ws2.Range("E26").Value = "Column 6"
$ws2.Range("E27").Value = "Column 6"
$ws2.Range("E28").Value = "Column 6"
$ws2.Range("E29").Value = "Column 6"
$ws2.Range("E30").Value = "Column 6"

# --- Update selection / active sheet to match the saved view state ---
$ws1.Range("E12").Select()

$ws2.Activate()
$ws2.Range("E31").Select()
